# Apply updated crypto price/volume values (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.037.17"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").Value = "1.563.65"
$ws.Range("E3").Value = "  +0.89%  "
$ws.Range("E4").Value = "  +0.36%  "
$ws.Range("D5").Value = "'208.51"
$ws.Range("E5").Value = "  +0.93%  "
$ws.Range("E6").Value = "  +0.65%  "
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("D8").Value = "'22.06"
$ws.Range("E8").Value = "  -0.72%  "
$ws.Range("E9").Value = "  +1.37%  "
$ws.Range("D10").Value = "'0.0598"
$ws.Range("E10").Value = "  +1.83%  "
$ws.Range("D12").Value = "1.785.16"
$ws.Range("E12").Value = "  +0.83%  "
$ws.Range("D13").Value = "1.562.32"
$ws.Range("E13").Value = "  +0.78%  "
$ws.Range("D14").Value = "'3.75"
$ws.Range("E14").Value = "  +0.40%  "
$ws.Range("D15").Value = "'0.520"
$ws.Range("E15").Value = "  +0.23%  "
$ws.Range("D16").Value = "27.024.95"
$ws.Range("E16").Value = "  +0.43%  "
$ws.Range("D17").Value = "'61.89"
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("E18").Value = "  +1.04%  "
$ws.Range("D19").Value = "'216.13"
$ws.Range("E19").Value = "  -0.71%  "
$ws.Range("D20").Value = "'7.38"
$ws.Range("E20").Value = "  +1.80%  "
$ws.Range("D21").Value = "'1.00"
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("D22").Value = "'4.15"
$ws.Range("E22").Value = "  +2.14%  "
$ws.Range("D23").Value = "'9.22"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  -0.51%  "
$ws.Range("D25").Value = "'153.78"
$ws.Range("E25").Value = "  -0.34%  "
$ws.Range("D26").Value = "'6.61"
$ws.Range("E26").Value = "  -0.23%  "
$ws.Range("D27").Value = "'15.06"
$ws.Range("E27").Value = "  +0.94%  "
$ws.Range("E28").Value = "  +1.64%  "
$ws.Range("D29").Value = "'1.01"
$ws.Range("E29").Value = "  +0.30%  "
$ws.Range("E30").Value = "  +1.53%  "
$ws.Range("E31").Value = "  +3.65%  "
$ws.Range("D32").Value = "'3.23"
$ws.Range("E32").Value = "  +0.10%  "
$ws.Range("E33").Value = "  +4.27%  "
$ws.Range("D34").Value = "1.425.50"
$ws.Range("E34").Value = "  +0.61%  "
$ws.Range("E35").Value = "  +11.72%  "
$ws.Range("D36").Value = "'1.61"
$ws.Range("E36").Value = "  +1.39%  "
$ws.Range("E37").Value = "  +2.49%  "
$ws.Range("E38").Value = "  +1.40%  "
$ws.Range("E39").Value = "  +1.05%  "
$ws.Range("D40").Value = "'0.811"
$ws.Range("E40").Value = "  +0.34%  "
$ws.Range("E41").Value = "  +1.44%  "
$ws.Range("D42").Value = "'1.01"
$ws.Range("E42").Value = "  +0.35%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("D45").Value = "'64.67"
$ws.Range("E45").Value = "  +0.27%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").Value = "1.698.00"
$ws.Range("E47").Value = "  +0.81%  "
$ws.Range("D48").Value = "'86.57"
$ws.Range("E48").Value = "  -1.11%  "
$ws.Range("D49").Value = "0.0₆0103"
$ws.Range("E49").Value = "  +1.65%  "
$ws.Range("D50").Value = "'0.0518"
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("D51").Value = "'0.0961"
$ws.Range("E51").Value = "  +0.11%  "
